# Edit script: updates batch1 sheet links/hyperlinks and replaces
# the "expanded" sheet's data table with the newly-scraped rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # batch1
$ws2 = $wb.Worksheets.Item(2)   # expanded

# ---------------------------------------------------------------
# Sheet "batch1": update the URL list in column B (rows 2-22)
# ---------------------------------------------------------------

# Row 2 keeps its URL but becomes a real hyperlink.
$ws1.Range("B2").Value = "https://www.american.edu/about/strategic-plan/2019/index.cfm"
$null = $ws1.Hyperlinks.Add($ws1.Range("B2"), $ws1.Range("B2").Value2)
$ws1.Range("B2").Style = "Hyperlink"

# Rows 3-8 get new plain (non-hyperlinked) URLs.
$ws1.Range("B3").Value = "https://www.american.edu/centers/antiracism/thinking-freedom-series.cfm"
$ws1.Range("B3").Style = "Normal"

$ws1.Range("B4").Value = "https://www.american.edu/about/strategic-plan/2019/team.cfm"
$ws1.Range("B4").Style = "Normal"

$ws1.Range("B5").Value = "https://www.american.edu/about/sustainability/campus-greening/waste.cfm"
$ws1.Range("B5").Style = "Normal"

$ws1.Range("B6").Value = "https://www.american.edu/about/sustainability/milestones.cfm"
$ws1.Range("B6").Style = "Normal"

$ws1.Range("B7").Value = "https://www.american.edu/about/sustainability/plastic-reduction.cfm"
$ws1.Range("B7").Style = "Normal"

$ws1.Range("B8").Value = "https://www.american.edu/about/sustainability/sustainability-competitions.cfm"
$ws1.Range("B8").Style = "Normal"

# Row 9 becomes a hyperlink to the adjudication page.
$ws1.Range("B9").Value = "https://www.american.edu/academics/integrity/adjudication.cfm"
$null = $ws1.Hyperlinks.Add($ws1.Range("B9"), $ws1.Range("B9").Value2)
$ws1.Range("B9").Style = "Hyperlink"

# Row 10 is a plain URL.
$ws1.Range("B10").Value = "https://www.american.edu/academics/integrity/code.cfm"
$ws1.Range("B10").Style = "Normal"

# Row 11 becomes a hyperlink to the faculty page.
$ws1.Range("B11").Value = "https://www.american.edu/academics/integrity/faculty.cfm"
$null = $ws1.Hyperlinks.Add($ws1.Range("B11"), $ws1.Range("B11").Value2)
$ws1.Range("B11").Style = "Hyperlink"

# Rows 12-22 no longer hold data; clear their contents (style stays s=2).
$ws1.Range("B12:B22").ClearContents()

# batch1 becomes the active/selected sheet, with B2:B11 selected.
$null = $ws1.Activate()
$null = $ws1.Range("B2:B11").Select()

# ---------------------------------------------------------------
# Sheet "expanded": replace the data table (rows 12 onward) with the
# refreshed set of URL / element-id / contains-DL-DT-Form-Table rows.
# ---------------------------------------------------------------

$u = "https://www.american.edu/centers/antiracism/thinking-freedom-series.cfm"
$ws2.Range("A12").Value = $u
$ws2.Range("B12").Value = "txt-5675573"
$ws2.Range("C12").Value = "❌"
$ws2.Range("A13").Value = $u
$ws2.Range("B13").Value = "txt-5674890"
$ws2.Range("C13").Value = "❌"
$ws2.Range("A14").Value = $u
$ws2.Range("B14").Value = "cs_control_5674895-collapse"
$ws2.Range("C14").Value = "✅"
$ws2.Range("A15").Value = $u
$ws2.Range("B15").Value = "cs_control_5674921-collapse"
$ws2.Range("C15").Value = "✅"
$ws2.Range("A16").Value = $u
$ws2.Range("B16").Value = "cs_control_5674922-collapse"
$ws2.Range("C16").Value = "✅"
$ws2.Range("A17").Value = $u
$ws2.Range("B17").Value = "cs_control_5674924-collapse"
$ws2.Range("C17").Value = "✅"
$ws2.Range("A18").Value = $u
$ws2.Range("B18").Value = "cs_control_5674958-collapse"
$ws2.Range("C18").Value = "✅"
$ws2.Range("A19").Value = $u
$ws2.Range("B19").Value = "txt-5691872"
$ws2.Range("C19").Value = "❌"
$u = "https://www.american.edu/about/strategic-plan/2019/team.cfm"
$ws2.Range("A20").Value = $u
$ws2.Range("B20").Value = "implementation-teams"
$ws2.Range("C20").Value = "❌"
$ws2.Range("A21").Value = $u
$ws2.Range("B21").Value = "txt-5433017"
$ws2.Range("C21").Value = "✅"
$ws2.Range("A22").Value = $u
$ws2.Range("B22").Value = "txt-5433068"
$ws2.Range("C22").Value = "❌"
$ws2.Range("A23").Value = $u
$ws2.Range("B23").Value = "txt-5433019"
$ws2.Range("C23").Value = "✅"
$ws2.Range("A24").Value = $u
$ws2.Range("B24").Value = "txt-5433021"
$ws2.Range("C24").Value = "✅"
$ws2.Range("A25").Value = $u
$ws2.Range("B25").Value = "txt-5433055"
$ws2.Range("C25").Value = "❌"
$ws2.Range("A26").Value = $u
$ws2.Range("B26").Value = "txt-5429100"
$ws2.Range("C26").Value = "✅"
$ws2.Range("A27").Value = $u
$ws2.Range("B27").Value = "txt-5429101"
$ws2.Range("C27").Value = "✅"
$u = "https://www.american.edu/about/sustainability/campus-greening/waste.cfm"
$ws2.Range("A28").Value = $u
$ws2.Range("B28").Value = "we-39-re-on-our-way-to-zero-waste"
$ws2.Range("C28").Value = "❌"
$ws2.Range("A29").Value = $u
$ws2.Range("B29").Value = "txt-4635685"
$ws2.Range("C29").Value = "❌"
$ws2.Range("A30").Value = $u
$ws2.Range("B30").Value = "txt-4704895"
$ws2.Range("C30").Value = "✅"
$ws2.Range("A31").Value = $u
$ws2.Range("B31").Value = "txt-4704907"
$ws2.Range("C31").Value = "❌"
$ws2.Range("A32").Value = $u
$ws2.Range("B32").Value = "txt-4781701"
$ws2.Range("C32").Value = "❌"
$ws2.Range("A33").Value = $u
$ws2.Range("B33").Value = "txt-4781709"
$ws2.Range("C33").Value = "❌"
$ws2.Range("A34").Value = $u
$ws2.Range("B34").Value = "txt-4781757"
$ws2.Range("C34").Value = "❌"
$ws2.Range("A35").Value = $u
$ws2.Range("B35").Value = "txt-4781805"
$ws2.Range("C35").Value = "❌"
$ws2.Range("A36").Value = $u
$ws2.Range("B36").Value = "txt-4704949"
$ws2.Range("C36").Value = "✅"
$ws2.Range("A37").Value = $u
$ws2.Range("B37").Value = "cs_control_5524586-collapse"
$ws2.Range("C37").Value = "❌"
$ws2.Range("A38").Value = $u
$ws2.Range("B38").Value = "cs_control_5524587-collapse"
$ws2.Range("C38").Value = "❌"
$ws2.Range("A39").Value = $u
$ws2.Range("B39").Value = "cs_control_5524591-collapse"
$ws2.Range("C39").Value = "✅"
$ws2.Range("A40").Value = $u
$ws2.Range("B40").Value = "txt-6877899"
$ws2.Range("C40").Value = "❌"
$u = "https://www.american.edu/about/sustainability/milestones.cfm"
$ws2.Range("A41").Value = $u
$ws2.Range("B41").Value = "txt-4613598"
$ws2.Range("C41").Value = "❌"
$ws2.Range("A42").Value = $u
$ws2.Range("B42").Value = "txt-6836132"
$ws2.Range("C42").Value = "❌"
$ws2.Range("A43").Value = $u
$ws2.Range("B43").Value = "txt-6836140"
$ws2.Range("C43").Value = "❌"
$ws2.Range("A44").Value = $u
$ws2.Range("B44").Value = "txt-6836112"
$ws2.Range("C44").Value = "❌"
$ws2.Range("A45").Value = $u
$ws2.Range("B45").Value = "txt-5734209"
$ws2.Range("C45").Value = "✅"
$ws2.Range("A46").Value = $u
$ws2.Range("B46").Value = "txt-5507286"
$ws2.Range("C46").Value = "✅"
$ws2.Range("A47").Value = $u
$ws2.Range("B47").Value = "txt-5507288"
$ws2.Range("C47").Value = "❌"
$ws2.Range("A48").Value = $u
$ws2.Range("B48").Value = "txt-5734729"
$ws2.Range("C48").Value = "❌"
$ws2.Range("A49").Value = $u
$ws2.Range("B49").Value = "txt-5507298"
$ws2.Range("C49").Value = "❌"
$ws2.Range("A50").Value = $u
$ws2.Range("B50").Value = "txt-5507289"
$ws2.Range("C50").Value = "✅"
$ws2.Range("A51").Value = $u
$ws2.Range("B51").Value = "txt-5507290"
$ws2.Range("C51").Value = "✅"
$ws2.Range("A52").Value = $u
$ws2.Range("B52").Value = "txt-5507299"
$ws2.Range("C52").Value = "❌"
$ws2.Range("A53").Value = $u
$ws2.Range("B53").Value = "txt-5507300"
$ws2.Range("C53").Value = "❌"
$ws2.Range("A54").Value = $u
$ws2.Range("B54").Value = "txt-5514819"
$ws2.Range("C54").Value = "❌"
$ws2.Range("A55").Value = $u
$ws2.Range("B55").Value = "txt-5515495"
$ws2.Range("C55").Value = "❌"
$ws2.Range("A56").Value = $u
$ws2.Range("B56").Value = "txt-5507320"
$ws2.Range("C56").Value = "❌"
$ws2.Range("A57").Value = $u
$ws2.Range("B57").Value = "txt-5507301"
$ws2.Range("C57").Value = "❌"
$ws2.Range("A58").Value = $u
$ws2.Range("B58").Value = "txt-5507422"
$ws2.Range("C58").Value = "❌"
$ws2.Range("A59").Value = $u
$ws2.Range("B59").Value = "txt-5515486"
$ws2.Range("C59").Value = "✅"
$ws2.Range("A60").Value = $u
$ws2.Range("B60").Value = "txt-5515488"
$ws2.Range("C60").Value = "❌"
$ws2.Range("A61").Value = $u
$ws2.Range("B61").Value = "txt-5515519"
$ws2.Range("C61").Value = "❌"
$ws2.Range("A62").Value = $u
$ws2.Range("B62").Value = "txt-5515510"
$ws2.Range("C62").Value = "❌"
$ws2.Range("A63").Value = $u
$ws2.Range("B63").Value = "txt-5507331"
$ws2.Range("C63").Value = "❌"
$ws2.Range("A64").Value = $u
$ws2.Range("B64").Value = "txt-5515472"
$ws2.Range("C64").Value = "❌"
$ws2.Range("A65").Value = $u
$ws2.Range("B65").Value = "txt-5507347"
$ws2.Range("C65").Value = "❌"
$ws2.Range("A66").Value = $u
$ws2.Range("B66").Value = "txt-5515525"
$ws2.Range("C66").Value = "❌"
$ws2.Range("A67").Value = $u
$ws2.Range("B67").Value = "txt-5515513"
$ws2.Range("C67").Value = "❌"
$ws2.Range("A68").Value = $u
$ws2.Range("B68").Value = "txt-5515465"
$ws2.Range("C68").Value = "❌"
$ws2.Range("A69").Value = $u
$ws2.Range("B69").Value = "txt-5507353"
$ws2.Range("C69").Value = "✅"
$ws2.Range("A70").Value = $u
$ws2.Range("B70").Value = "txt-5507380"
$ws2.Range("C70").Value = "✅"
$ws2.Range("A71").Value = $u
$ws2.Range("B71").Value = "txt-5515532"
$ws2.Range("C71").Value = "❌"
$u = "https://www.american.edu/about/sustainability/plastic-reduction.cfm"
$ws2.Range("A72").Value = $u
$ws2.Range("B72").Value = "benefits-of-reusable-water-bottles"
$ws2.Range("C72").Value = "❌"
$ws2.Range("A73").Value = $u
$ws2.Range("B73").Value = "dc-water"
$ws2.Range("C73").Value = "❌"
$ws2.Range("A74").Value = $u
$ws2.Range("B74").Value = "cs_control_6346253-collapse"
$ws2.Range("C74").Value = "✅"
$ws2.Range("A75").Value = $u
$ws2.Range("B75").Value = "cs_control_6361607-collapse"
$ws2.Range("C75").Value = "✅"
$ws2.Range("A76").Value = $u
$ws2.Range("B76").Value = "cs_control_6346259-collapse"
$ws2.Range("C76").Value = "✅"
$ws2.Range("A77").Value = $u
$ws2.Range("B77").Value = "cs_control_6346261-collapse"
$ws2.Range("C77").Value = "✅"
$ws2.Range("A78").Value = $u
$ws2.Range("B78").Value = "cs_control_6361592-collapse"
$ws2.Range("C78").Value = "✅"
$ws2.Range("A79").Value = $u
$ws2.Range("B79").Value = "cs_control_6346274-collapse"
$ws2.Range("C79").Value = "✅"
$ws2.Range("A80").Value = $u
$ws2.Range("B80").Value = "cs_control_6346263-collapse"
$ws2.Range("C80").Value = "✅"
$ws2.Range("A81").Value = $u
$ws2.Range("B81").Value = "cs_control_6361600-collapse"
$ws2.Range("C81").Value = "✅"
$ws2.Range("A82").Value = $u
$ws2.Range("B82").Value = "cs_control_6346271-collapse"
$ws2.Range("C82").Value = "✅"
$ws2.Range("A83").Value = $u
$ws2.Range("B83").Value = "cs_control_6346279-collapse"
$ws2.Range("C83").Value = "✅"
$ws2.Range("A84").Value = $u
$ws2.Range("B84").Value = "cs_control_6361613-collapse"
$ws2.Range("C84").Value = "✅"
$ws2.Range("A85").Value = $u
$ws2.Range("B85").Value = "cs_control_6346252-collapse"
$ws2.Range("C85").Value = "✅"
$ws2.Range("A86").Value = $u
$ws2.Range("B86").Value = "cs_control_6346257-collapse"
$ws2.Range("C86").Value = "✅"
$ws2.Range("A87").Value = $u
$ws2.Range("B87").Value = "cs_control_6346269-collapse"
$ws2.Range("C87").Value = "✅"
$u = "https://www.american.edu/about/sustainability/sustainability-competitions.cfm"
$ws2.Range("A88").Value = $u
$ws2.Range("B88").Value = "campus-race-to-zero-waste"
$ws2.Range("C88").Value = "❌"
$ws2.Range("A89").Value = $u
$ws2.Range("B89").Value = "cs_control_5896649-collapse"
$ws2.Range("C89").Value = "✅"
$ws2.Range("A90").Value = $u
$ws2.Range("B90").Value = "cs_control_5896651-collapse"
$ws2.Range("C90").Value = "❌"
$ws2.Range("A91").Value = $u
$ws2.Range("B91").Value = "cs_control_5948475-collapse"
$ws2.Range("C91").Value = "❌"
$u = "https://www.american.edu/academics/integrity/adjudication.cfm"
$ws2.Range("A92").Value = $u
$ws2.Range("B92").Value = "adjudication-process"
$ws2.Range("C92").Value = "✅"
$u = "https://www.american.edu/academics/integrity/code.cfm"
$ws2.Range("A93").Value = $u
$ws2.Range("B93").Value = "academic-integrity-code"
$ws2.Range("C93").Value = "❌"
$ws2.Range("A94").Value = $u
$ws2.Range("B94").Value = "watch-academic-integrity-videos"
$ws2.Range("C94").Value = "✅"
$ws2.Range("A95").Value = $u
$ws2.Range("B95").Value = "txt-5644841"
$ws2.Range("C95").Value = "✅"
$ws2.Range("A96").Value = $u
$ws2.Range("B96").Value = "txt-5644844"
$ws2.Range("C96").Value = "✅"
$u = "https://www.american.edu/academics/integrity/faculty.cfm"
$ws2.Range("A97").Value = $u
$ws2.Range("B97").Value = "faculty-resources"
$ws2.Range("C97").Value = "✅"

Write-Host "edit complete"
